# Fruta / hortaliza, semanal
#
# Insert two new weekly price-report rows (new row 14 "Primera" and new
# row 15 "Segunda", both dated 2021-12-03) into the "Macroferia Regional
# de Talca - Arandano (blue)" sheet, pushing the existing rows 14-36 down
# to rows 16-38.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows above the current row 14 (shifts old rows 14.. down by 2).
$ws.Range("A14:A15").EntireRow.Insert()

# New row 14: "Primera" quality entry for 2021-12-03.
$ws.Cells.Item(14, 1).Value = 5
$ws.Cells.Item(14, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(14, 3).Value = "Maule"
$ws.Cells.Item(14, 4).Value = "2021-12-03"
$ws.Cells.Item(14, 5).Value = 7
$ws.Cells.Item(14, 6).Value = "Fruta"
$ws.Cells.Item(14, 7).Value = 100101
$ws.Cells.Item(14, 8).Value = "Berries"
$ws.Cells.Item(14, 9).Value = 100101001
$ws.Cells.Item(14, 10).Value = "Arándano (blue)"
$ws.Cells.Item(14, 11).Value = "Sin especificar"
$ws.Cells.Item(14, 12).Value = "Primera"
$ws.Cells.Item(14, 13).Value = 180
$ws.Cells.Item(14, 14).Value = 3600
$ws.Cells.Item(14, 15).Value = 3600
$ws.Cells.Item(14, 16).Value = 3600
$ws.Cells.Item(14, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(14, 18).Value = "Provincia de Linares"
$ws.Cells.Item(14, 19).Value = 1800
$ws.Cells.Item(14, 20).Value = 2

# New row 15: "Segunda" quality entry for 2021-12-03.
$ws.Cells.Item(15, 1).Value = 5
$ws.Cells.Item(15, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(15, 3).Value = "Maule"
$ws.Cells.Item(15, 4).Value = "2021-12-03"
$ws.Cells.Item(15, 5).Value = 7
$ws.Cells.Item(15, 6).Value = "Fruta"
$ws.Cells.Item(15, 7).Value = 100101
$ws.Cells.Item(15, 8).Value = "Berries"
$ws.Cells.Item(15, 9).Value = 100101001
$ws.Cells.Item(15, 10).Value = "Arándano (blue)"
$ws.Cells.Item(15, 11).Value = "Sin especificar"
$ws.Cells.Item(15, 12).Value = "Segunda"
$ws.Cells.Item(15, 13).Value = 100
$ws.Cells.Item(15, 14).Value = 3000
$ws.Cells.Item(15, 15).Value = 3000
$ws.Cells.Item(15, 16).Value = 3000
$ws.Cells.Item(15, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(15, 18).Value = "Provincia de Linares"
$ws.Cells.Item(15, 19).Value = 1500
$ws.Cells.Item(15, 20).Value = 2
